# Updates computed-profit figures (currentAveragePrice / profit columns H-N)
# across the per-job leve-profit sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR),
# reflecting refreshed market-board pricing data from the scheduled data-pull
# runner. Values below are written directly onto the affected cells.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 11689
$ws.Range("I47").Value = 16033.5
$ws.Range("J47").Value = 3000
$ws.Range("K47").Value = 16033.5
$ws.Range("L47").Value = 3000
$ws.Range("M47").Value = -15061.5
$ws.Range("N47").Value = -4944
$ws.Range("H53").Value = 573.03845
$ws.Range("I53").Value = 171.08333
$ws.Range("J53").Value = 917.5714
$ws.Range("K53").Value = 171.08333
$ws.Range("L53").Value = 917.5714
$ws.Range("M53").Value = 465.91667
$ws.Range("N53").Value = -2191.5714
$ws.Range("H55").Value = 344.2857
$ws.Range("I55").Value = 385
$ws.Range("J55").Value = 100
$ws.Range("K55").Value = 385
$ws.Range("L55").Value = 100
$ws.Range("M55").Value = -171
$ws.Range("N55").Value = -528
$ws.Range("H132").Value = 15275.3125
$ws.Range("I132").Value = 9627
$ws.Range("K132").Value = 28881
$ws.Range("M132").Value = -26351

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9202.27
$ws.Range("I32").Value = 9489.666999999999
$ws.Range("J32").Value = 5753.5
$ws.Range("K32").Value = 9489.666999999999
$ws.Range("L32").Value = 5753.5
$ws.Range("M32").Value = -9202.666999999999
$ws.Range("N32").Value = -6327.5
$ws.Range("H88").Value = 2625
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 2833.3333
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 2833.3333
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -3645.3333
$ws.Range("H91").Value = 2625
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 2833.3333
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 2833.3333
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -5641.3333
$ws.Range("H101").Value = 18275.5
$ws.Range("J101").Value = 18275.5
$ws.Range("L101").Value = 18275.5
$ws.Range("N101").Value = -24765.5
$ws.Range("H109").Value = 26000
$ws.Range("J109").Value = 26000
$ws.Range("L109").Value = 26000
$ws.Range("N109").Value = -28774
$ws.Range("H122").Value = 2267.3928
$ws.Range("I122").Value = 1466.2174
$ws.Range("J122").Value = 5952.8
$ws.Range("K122").Value = 4398.6522
$ws.Range("L122").Value = 17858.4
$ws.Range("M122").Value = -1948.6522
$ws.Range("N122").Value = -22758.4
$ws.Range("H132").Value = 4144.041
$ws.Range("I132").Value = 1747.36
$ws.Range("J132").Value = 6640.5835
$ws.Range("K132").Value = 5242.08
$ws.Range("L132").Value = 19921.7505
$ws.Range("M132").Value = -2712.08
$ws.Range("N132").Value = -24981.7505

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2346.0625
$ws.Range("I86").Value = 2007.4
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 2007.4
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -884.4000000000001
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 2346.0625
$ws.Range("I89").Value = 2007.4
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 10037
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -4421
$ws.Range("N89").Value = -23732
$ws.Range("H94").Value = 887.1053000000001
$ws.Range("I94").Value = 755.3570999999999
$ws.Range("J94").Value = 1256
$ws.Range("K94").Value = 755.3570999999999
$ws.Range("L94").Value = 1256
$ws.Range("M94").Value = -304.3570999999999
$ws.Range("N94").Value = -2158
$ws.Range("H100").Value = 22666.666
$ws.Range("J100").Value = 22666.666
$ws.Range("L100").Value = 22666.666
$ws.Range("N100").Value = -24830.666
$ws.Range("H134").Value = 5250.125
$ws.Range("I134").Value = 2924.4119
$ws.Range("J134").Value = 6969.1304
$ws.Range("K134").Value = 8773.235700000001
$ws.Range("L134").Value = 20907.3912
$ws.Range("M134").Value = -6238.235700000001
$ws.Range("N134").Value = -25977.3912

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2529.6445
$ws.Range("I31").Value = 1758.8276
$ws.Range("J31").Value = 3926.75
$ws.Range("K31").Value = 1758.8276
$ws.Range("L31").Value = 3926.75
$ws.Range("M31").Value = -1463.8276
$ws.Range("N31").Value = -4516.75
$ws.Range("H34").Value = 2529.6445
$ws.Range("I34").Value = 1758.8276
$ws.Range("J34").Value = 3926.75
$ws.Range("K34").Value = 1758.8276
$ws.Range("L34").Value = 3926.75
$ws.Range("M34").Value = -1556.8276
$ws.Range("N34").Value = -4330.75
$ws.Range("H58").Value = 1928763.9
$ws.Range("I58").Value = 2125.5
$ws.Range("J58").Value = 8350891.5
$ws.Range("K58").Value = 2125.5
$ws.Range("L58").Value = 8350891.5
$ws.Range("M58").Value = -1922.5
$ws.Range("N58").Value = -8351297.5
$ws.Range("H107").Value = 1830.5
$ws.Range("J107").Value = 3247.5
$ws.Range("L107").Value = 3247.5
$ws.Range("N107").Value = -7087.5
$ws.Range("H136").Value = 1928763.9
$ws.Range("I136").Value = 2125.5
$ws.Range("J136").Value = 8350891.5
$ws.Range("K136").Value = 6376.5
$ws.Range("L136").Value = 25052674.5
$ws.Range("M136").Value = -3826.5
$ws.Range("N136").Value = -25057774.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1159.9412
$ws.Range("I131").Value = 331.81818
$ws.Range("J131").Value = 2678.1667
$ws.Range("K131").Value = 995.45454
$ws.Range("L131").Value = 8034.500100000001
$ws.Range("M131").Value = 4044.54546
$ws.Range("N131").Value = -18114.5001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 5000
$ws.Range("J44").Value = 5000
$ws.Range("L44").Value = 5000
$ws.Range("N44").Value = -5912
$ws.Range("H68").Value = 2448.3333
$ws.Range("I68").Value = 2153.3333
$ws.Range("J68").Value = 3333.3333
$ws.Range("K68").Value = 2153.3333
$ws.Range("L68").Value = 3333.3333
$ws.Range("M68").Value = -1404.3333
$ws.Range("N68").Value = -4831.3333
$ws.Range("H71").Value = 2448.3333
$ws.Range("I71").Value = 2153.3333
$ws.Range("J71").Value = 3333.3333
$ws.Range("K71").Value = 10766.6665
$ws.Range("L71").Value = 16666.6665
$ws.Range("M71").Value = -7022.666499999999
$ws.Range("N71").Value = -24154.6665

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 16000
$ws.Range("J94").Value = 16000
$ws.Range("L94").Value = 16000
$ws.Range("N94").Value = -17802
$ws.Range("H96").Value = 5056.364
$ws.Range("I96").Value = 2323.077
$ws.Range("J96").Value = 9004.444
$ws.Range("K96").Value = 2323.077
$ws.Range("L96").Value = 9004.444
$ws.Range("M96").Value = -950.0770000000002
$ws.Range("N96").Value = -11750.444
$ws.Range("H97").Value = 10500
$ws.Range("J97").Value = 10500
$ws.Range("L97").Value = 10500
$ws.Range("N97").Value = -12482
$ws.Range("H101").Value = 8214.571
$ws.Range("J101").Value = 8214.571
$ws.Range("L101").Value = 8214.571
$ws.Range("N101").Value = -14704.571
$ws.Range("H132").Value = 1907.1538
$ws.Range("I132").Value = 1530.9546
$ws.Range("K132").Value = 4592.8638
$ws.Range("M132").Value = -2062.8638
$ws.Range("H136").Value = 7144511
$ws.Range("I136").Value = 22728372
$ws.Range("J136").Value = 1908.125
$ws.Range("K136").Value = 68185116
$ws.Range("L136").Value = 5724.375
$ws.Range("M136").Value = -68182566
$ws.Range("N136").Value = -10824.375
